# Apply crypto price/volume updates scraped on Tue Jun 20 23:01:34 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.212.95"
$ws.Range("E2").Value = "  +5.62%  "
$ws.Range("D3").Value = "1.785.02"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2679"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").Value = "1.781.05"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07047"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6284"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.665"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "79.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").Value = "28.194.13"
$ws.Range("E16").Value = "  +6.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9992"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007237"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.36%  "
$ws.Range("D21").Value = "2.009.20"
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.555"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.759"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.255"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("E28").Value = "  +2.98%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.187"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08275"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.768"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04903"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.087"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.614"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6529"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9461"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  +8.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.054"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.932"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01551"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3989"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("E45").Value = "  +3.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1215"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.75%  "
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.030"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.294"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.92%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("E51").Value = "  +2.17%  "
